$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D values are numeric-looking text (e.g. "1.001", "30.018.86") that must
# stay stored as text, matching the original inlineStr cells. Force text format
# before assigning, then restore the default "Normal" style so no stray number
# format / style index is left behind on the cell.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "30.018.86"
$ws.Range("D2").Style = "Normal"
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.902.76"
$ws.Range("D3").Style = "Normal"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7452"
$ws.Range("D5").Style = "Normal"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.22"
$ws.Range("D6").Style = "Normal"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3075"
$ws.Range("D8").Style = "Normal"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.68"
$ws.Range("D9").Style = "Normal"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.06912"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08050"
$ws.Range("D11").Style = "Normal"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7565"
$ws.Range("D12").Style = "Normal"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.903.68"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.243"
$ws.Range("D14").Style = "Normal"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "91.32"
$ws.Range("D15").Style = "Normal"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.190"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "30.014.38"
$ws.Range("D17").Style = "Normal"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.000007770"
$ws.Range("D19").Style = "Normal"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "237.00"
$ws.Range("D20").Style = "Normal"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "1.001"
$ws.Range("D21").Style = "Normal"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "2.155.19"
$ws.Range("D22").Style = "Normal"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "1.001"
$ws.Range("D23").Style = "Normal"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "167.99"
$ws.Range("D25").Style = "Normal"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.338"
$ws.Range("D26").Style = "Normal"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.1277"
$ws.Range("D28").Style = "Normal"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.352"
$ws.Range("D30").Style = "Normal"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.533"
$ws.Range("D31").Style = "Normal"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.308"
$ws.Range("D32").Style = "Normal"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "4.047"
$ws.Range("D33").Style = "Normal"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.05265"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.284"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7405"
$ws.Range("D36").Style = "Normal"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.726"
$ws.Range("D37").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "2.765"
$ws.Range("D39").Style = "Normal"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.256"
$ws.Range("D40").Style = "Normal"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4466"
$ws.Range("D41").Style = "Normal"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.66"
$ws.Range("D42").Style = "Normal"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.949"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "1.001"
$ws.Range("D44").Style = "Normal"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.8332"
$ws.Range("D45").Style = "Normal"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "7.737"
$ws.Range("D46").Style = "Normal"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "101.49"
$ws.Range("D47").Style = "Normal"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "2.058.32"
$ws.Range("D49").Style = "Normal"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "36.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.05987"
$ws.Range("D51").Style = "Normal"

# Remaining text cells (coin name/link in row 51, and the Volume(1h) percentages)
# are safe to assign directly -- they contain non-numeric characters (%, spaces,
# letters, slashes) so Excel keeps them as text without any extra coercion.
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("E3").Value = "  -1.41%  "
$ws.Range("E4").Value = "  -0.27%  "
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  -0.76%  "
$ws.Range("E7").Value = "  -0.24%  "
$ws.Range("E8").Value = "  -2.46%  "
$ws.Range("E9").Value = "  -6.35%  "
$ws.Range("E10").Value = "  -1.53%  "
$ws.Range("E11").Value = "  +0.19%  "
$ws.Range("E12").Value = "  -2.54%  "
$ws.Range("E14").Value = "  -2.51%  "
$ws.Range("E15").Value = "  -1.96%  "
$ws.Range("E16").Value = "  +3.30%  "
$ws.Range("E17").Value = "  -1.07%  "
$ws.Range("E18").Value = "  -3.08%  "
$ws.Range("E19").Value = "  -1.76%  "
$ws.Range("E20").Value = "  -5.49%  "
$ws.Range("E21").Value = "  -0.19%  "
$ws.Range("E22").Value = "  -0.67%  "
$ws.Range("E23").Value = "  -0.20%  "
$ws.Range("E24").Value = "  +7.10%  "
$ws.Range("E25").Value = "  +1.57%  "
$ws.Range("E26").Value = "  -2.27%  "
$ws.Range("E27").Value = "  -0.97%  "
$ws.Range("E28").Value = "  -0.99%  "
$ws.Range("E29").Value = "  -5.15%  "
$ws.Range("E30").Value = "  -0.38%  "
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("E32").Value = "  -2.13%  "
$ws.Range("E33").Value = "  -1.44%  "
$ws.Range("E34").Value = "  +0.96%  "
$ws.Range("E35").Value = "  -1.96%  "
$ws.Range("E36").Value = "  -1.79%  "
$ws.Range("E37").Value = "  -1.53%  "
$ws.Range("E38").Value = "  -0.17%  "
$ws.Range("E39").Value = "  -0.99%  "
$ws.Range("E40").Value = "  -3.79%  "
$ws.Range("E41").Value = "  -0.64%  "
$ws.Range("E42").Value = "  -5.03%  "
$ws.Range("E44").Value = "  -0.19%  "
$ws.Range("E45").Value = "  -1.08%  "
$ws.Range("E46").Value = "  +0.95%  "
$ws.Range("E47").Value = "  +0.22%  "
$ws.Range("E48").Value = "  -1.77%  "
$ws.Range("E49").Value = "  -1.05%  "
$ws.Range("E50").Value = "  -2.76%  "
$ws.Range("B51").Value = "Cronos"
$ws.Range("C51").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("E51").Value = "  -0.48%  "
